# Updated 2D training schedules, no break screen
# Adds a new "break_on_off" column (L) to Sheet1, flagging the trials
# right before a break screen (rows 19, 37 and 54 -> trial numbers 18, 36, 53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("L1").Value = "break_on_off"

# Values for rows 2..73 (trials 1..72); 1 marks a break, 0 otherwise.
$breakOnOff = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $breakOnOff.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $breakOnOff[$i]
}
